# Rotate the contents of columns D,E,F,G (codeforiati:group-name,
# codeforiati:category-name, codeforiati:group-code, codeforiati:category-code)
# for every row so that:
#   new D = old G   (codeforiati:category-code)
#   new E = old F   (codeforiati:group-code)
#   new F = old D   (codeforiati:group-name)
#   new G = old E   (codeforiati:category-name)
#
# This also re-orders the header row (row 1) the same way, since the header
# labels themselves move along with the data columns.
#
# Cells D:G are plain shared-string ("text") cells, even though many of them
# look numeric (e.g. "110"). Using Range.Value with literal strings causes
# Excel to re-interpret numeric-looking text as a Double, which would change
# the stored cell type. Range.Copy preserves the original cell's value type
# (and lack of any special formatting) exactly, so it is used here together
# with a single scratch column (column Z) as a one-cell swap buffer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

$tempCol = 26  # column Z, used transiently as a swap buffer

for ($r = 1; $r -le $lastRow; $r++) {
    $cD = $ws.Cells.Item($r, 4)
    $cE = $ws.Cells.Item($r, 5)
    $cF = $ws.Cells.Item($r, 6)
    $cG = $ws.Cells.Item($r, 7)
    $cTemp = $ws.Cells.Item($r, $tempCol)

    $cD.Copy($cTemp)   # temp = D (old)
    $cG.Copy($cD)      # D = G (old)
    $cE.Copy($cG)      # G = E (old)
    $cF.Copy($cE)      # E = F (old)
    $cTemp.Copy($cF)   # F = temp (old D)

    $cTemp.ClearContents()
}
